$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: E1/F1 swap ---
$ws.Range("E1").Value = "UpdateNoofClasses"
$ws.Range("F1").Value = "UpdateBatchDescription"

# --- Row 2 updates ---
$ws.Range("G2").Value = "!@3*&Invalid"

# --- New column K value/header (set value order matches original authoring order) ---
$ws.Range("K2").Value = "This field should start with an alphabet and min 2 character."
$ws.Range("K2").Style = "Hyperlink"
$ws.Range("K2").Font.ThemeFont = 1
$ws.Range("K2").Font.ThemeColor = 1
$ws.Range("K2").WrapText = $true
$ws.Columns.Item(11).ColumnWidth = 17.17

$ws.Range("K1").Value = "errormsg"

# --- Row 3 new cell ---
$ws.Range("G3").Value = "     "

# --- Row 2 height (auto after wrap) ---
$ws.Rows.Item(2).RowHeight = 43.2

# --- Add new data into existing row 4 (A4,B4,D4) ---
$ws.Range("A4").Value = "Invalid without optionaldesc"
$ws.Range("B4").Value = "BatchTribe"
$ws.Range("D4").Value = 3

$ws.Range("A4").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("I4").WrapText = $true
$ws.Range("J4").WrapText = $true

$ws.Rows.Item(4).RowHeight = 40.8

# --- selection ---
$ws.Range("A4:XFD4").Select()
